# Generate Report for Handoff
# Updates Priority ("low" -> "ht") and Latest Handoff Datetime for the
# rows that were just handed off, on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7, Priority (E) low -> ht, Latest Handoff Datetime (H) updated
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-31 06:35:28"
}

# de-de sheet: rows 4-7, Priority (E) low -> ht, Latest Handoff Datetime (H) updated
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-31 06:35:33"
}
